# Update Fbln1-Itgb1 LR-pair TPM results: recomputed metrics for the
# existing Sending-cluster x Target-cluster combinations (ECs/FAPs/MuSCs),
# plus four new rows for the "Resolving-Mac" sending cluster (rows 14-17).
$rows = @(
    @("ECs","Fbln1","Itgb1","ECs",3,1,0.1919306666666667,0.575792,0.004053148510572174,0.004053148510572173,3,1,168.1098273333333,504.329482,0.2984182258032519,0.298418225803252,32.26543123330489,290.388881099744,0.001209533387442041,0.001209533387442041),
    @("ECs","Fbln1","Itgb1","FAPs",3,1,0.1919306666666667,0.575792,0.004053148510572174,0.004053148510572173,3,1,163.0062356666667,489.018707,0.2893586437755394,0.2893586437755394,31.28589548232711,281.573059340944,0.001172813556040012,0.001172813556040012),
    @("ECs","Fbln1","Itgb1","MuSCs",3,1,0.1919306666666667,0.575792,0.004053148510572174,0.004053148510572173,3,1,165.99353,497.98059,0.294661504941043,0.294661504941043,31.85924887525333,286.73323987728,0.001194306839874744,0.001194306839874744),
    @("ECs","Fbln1","Itgb1","Resolving-Mac",3,1,0.1919306666666667,0.575792,0.004053148510572174,0.004053148510572173,3,1,66.22673433333334,198.680203,0.1175616254801657,0.1175616254801657,12.71094127175289,114.398471445776,0.0004764947272153772,0.0004764947272153771),
    @("FAPs","Fbln1","Itgb1","ECs",3,1,43.24729533333333,129.741886,0.9132866243360881,0.9132866243360879,3,1,168.1098273333333,504.329482,0.2984182258032519,0.298418225803252,7270.295351120339,65432.65816008305,0.2725413740842165,0.2725413740842165),
    @("FAPs","Fbln1","Itgb1","FAPs",3,1,43.24729533333333,129.741886,0.9132866243360881,0.9132866243360879,3,1,163.0062356666667,489.018707,0.2893586437755394,0.2893586437755394,7049.578815051267,63446.2093354614,0.264267378996231,0.264267378996231),
    @("FAPs","Fbln1","Itgb1","MuSCs",3,1,43.24729533333333,129.741886,0.9132866243360881,0.9132866243360879,3,1,165.99353,497.98059,0.294661504941043,0.294661504941043,7178.771215332526,64608.94093799274,0.2691104111693967,0.2691104111693966),
    @("FAPs","Fbln1","Itgb1","Resolving-Mac",3,1,43.24729533333333,129.741886,0.9132866243360881,0.9132866243360879,3,1,66.22673433333334,198.680203,0.1175616254801657,0.1175616254801657,2864.127138675874,25777.14424808286,0.1073674600862439,0.1073674600862439),
    @("MuSCs","Fbln1","Itgb1","ECs",3,1,3.905830333333333,11.717491,0.0824824436502988,0.08248244365029879,3,1,168.1098273333333,504.329482,0.2984182258032519,0.298418225803252,656.6084629299623,5909.476166369661,0.02461426449403887,0.02461426449403887),
    @("MuSCs","Fbln1","Itgb1","FAPs",3,1,3.905830333333333,11.717491,0.0824824436502988,0.08248244365029879,3,1,163.0062356666667,489.018707,0.2893586437755394,0.2893586437755394,636.6746997893484,5730.072298104136,0.02386700802994282,0.02386700802994281),
    @("MuSCs","Fbln1","Itgb1","MuSCs",3,1,3.905830333333333,11.717491,0.0824824436502988,0.08248244365029879,3,1,165.99353,497.98059,0.294661504941043,0.294661504941043,648.3425646110766,5835.08308149969,0.02430440097721182,0.02430440097721182),
    @("MuSCs","Fbln1","Itgb1","Resolving-Mac",3,1,3.905830333333333,11.717491,0.0824824436502988,0.08248244365029879,3,1,66.22673433333334,198.680203,0.1175616254801657,0.1175616254801657,258.6703878367414,2328.033490530673,0.009696770149105296,0.009696770149105294),
    @("Resolving-Mac","Fbln1","Itgb1","ECs",1,0.3333333333333333,0.008418666666666666,0.025256,0.0001777835030410475,0.0001777835030410475,3,1,168.1098273333333,504.329482,0.2984182258032519,0.298418225803252,1.415260599710222,12.737345397392,0.00005305383755459644,0.00005305383755459645),
    @("Resolving-Mac","Fbln1","Itgb1","FAPs",1,0.3333333333333333,0.008418666666666666,0.025256,0.0001777835030410475,0.0001777835030410475,3,1,163.0062356666667,489.018707,0.2893586437755394,0.2893586437755394,1.372295162665778,12.350656463992,0.00005144319332562199,0.00005144319332562199),
    @("Resolving-Mac","Fbln1","Itgb1","MuSCs",1,0.3333333333333333,0.008418666666666666,0.025256,0.0001777835030410475,0.0001777835030410475,3,1,165.99353,497.98059,0.294661504941043,0.294661504941043,1.397444197893333,12.57699778104,0.00005238595455976555,0.00005238595455976555),
    @("Resolving-Mac","Fbln1","Itgb1","Resolving-Mac",1,0.3333333333333333,0.008418666666666666,0.025256,0.0001777835030410475,0.0001777835030410475,3,1,66.22673433333334,198.680203,0.1175616254801657,0.1175616254801657,0.5575408007742223,5.017867206968,0.00002090051760106352,0.00002090051760106352)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

